$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F-column (time_taken) timestamps on the "data" sheet ---
$newTimes = @(
  "2021-10-05 14:34:27.035364",
  "2021-10-05 14:34:27.035372",
  "2021-10-05 14:34:27.035375",
  "2021-10-05 14:34:27.035378",
  "2021-10-05 14:34:27.035381",
  "2021-10-05 14:34:27.035384",
  "2021-10-05 14:34:27.035387",
  "2021-10-05 14:34:27.035389",
  "2021-10-05 14:34:27.035392",
  "2021-10-05 14:34:27.035395",
  "2021-10-05 14:34:27.035398",
  "2021-10-05 14:34:27.035401",
  "2021-10-05 14:34:27.035403",
  "2021-10-05 14:34:27.035406",
  "2021-10-05 14:34:27.035409",
  "2021-10-05 14:34:27.035411",
  "2021-10-05 14:34:27.035414",
  "2021-10-05 14:34:27.035417",
  "2021-10-05 14:34:27.035419",
  "2021-10-05 14:34:27.035422",
  "2021-10-05 14:34:27.035425",
  "2021-10-05 14:34:27.035427",
  "2021-10-05 14:34:27.035430",
  "2021-10-05 14:34:27.035432",
  "2021-10-05 14:34:27.035436",
  "2021-10-05 14:34:27.035438",
  "2021-10-05 14:34:27.035441",
  "2021-10-05 14:34:27.035443",
  "2021-10-05 14:34:27.035446",
  "2021-10-05 14:34:27.035449",
  "2021-10-05 14:34:27.035451",
  "2021-10-05 14:34:27.035454",
  "2021-10-05 14:34:27.035457",
  "2021-10-05 14:34:27.035460",
  "2021-10-05 14:34:27.035462",
  "2021-10-05 14:34:27.035465",
  "2021-10-05 14:34:27.035468",
  "2021-10-05 14:34:27.035470",
  "2021-10-05 14:34:27.035473",
  "2021-10-05 14:34:27.035476",
  "2021-10-05 14:34:27.035479",
  "2021-10-05 14:34:27.035482",
  "2021-10-05 14:34:27.035485",
  "2021-10-05 14:34:27.035487",
  "2021-10-05 14:34:27.035490",
  "2021-10-05 14:34:27.035493",
  "2021-10-05 14:34:27.035495",
  "2021-10-05 14:34:27.035498"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" sheet, placed after "data" ---
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"
$meta.Outline.SummaryRow = 1
$meta.Outline.SummaryColumn = 1

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Lymphoedema_syndromic"
$meta.Range("C2").Value = 3098
# data_version ("0.11") is a version label, not a number -> force text storage,
# then restore the default (unstyled) look by pasting the plain formatting
# from a normal, unstyled cell on the "data" sheet over it.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.11"
$dataSheet.Range("B2").Copy()
$meta.Range("D2").PasteSpecial(-4122)
$meta.Range("E2").Value = "2021-09-22T23:54:32.996455Z"
$meta.Range("F2").Value = "2021-10-05 14:34:27.031436"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3098/?format=json"

# Match the header-row / index-column styling used on the "data" sheet (bold,
# centered, bordered -> reuses the existing style index instead of minting a
# new one).
$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$dataSheet.Activate()
$dataSheet.Range("A1").Select()
